# Add "当月放款" / "当月放款笔数" columns to the "对公产品台账202404" sheet
# (columns L, M) and to the "个人经营贷202404" sheet (columns N, O).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 4: 对公产品台账202404  -> new columns L (当月放款), M (当月放款笔数)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Copy the header style (bold / centered / bordered) from the last existing
# header cell (K1) onto the two new header cells before writing their text.
$ws4.Range("K1").Copy()
$ws4.Range("L1:M1").PasteSpecial(-4122)  # xlPasteFormats

$ws4.Cells.Item(1, 12).Value = "当月放款"
$ws4.Cells.Item(1, 13).Value = "当月放款笔数"

$ws4Data = @{
    2  = @(206.9, 1)
    3  = @(3270.47, 71)
    4  = @(100, 1)
    5  = @(170, 1)
    6  = @(0, 0)
    7  = @(13310.25, 405)
    8  = @(0, 0)
    9  = @(0, 0)
    10 = @(0, 0)
    11 = @(0, 0)
    12 = @(0, 0)
    13 = @(0, 0)
    14 = @(2538.7, 23)
    15 = @(195.5, 1)
    16 = @(19791.82, 503)
}

foreach ($row in $ws4Data.Keys) {
    $pair = $ws4Data[$row]
    $ws4.Cells.Item($row, 12).Value = $pair[0]
    $ws4.Cells.Item($row, 13).Value = $pair[1]
}

# ---------------------------------------------------------------------
# Sheet 5: 个人经营贷202404  -> new columns N (当月放款), O (当月放款笔数)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("M1").Copy()
$ws5.Range("N1:O1").PasteSpecial(-4122)  # xlPasteFormats

$ws5.Cells.Item(1, 14).Value = "当月放款"
$ws5.Cells.Item(1, 15).Value = "当月放款笔数"

$ws5Data = @{
    2  = @(0, 0)
    3  = @(19732, 102)
    4  = @(18931, 77)
    5  = @(139.65, 4)
    6  = @(1954.6, 37)
    7  = @(1574.06, 81)
    8  = @(335, 5)
    9  = @(5942.7, 368)
    10 = @(6618.5, 356)
    11 = @(0, 0)
    12 = @(299, 3)
    13 = @(1160, 6)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @(56686.50999999999, 1039)
}

foreach ($row in $ws5Data.Keys) {
    $pair = $ws5Data[$row]
    $ws5.Cells.Item($row, 14).Value = $pair[0]
    $ws5.Cells.Item($row, 15).Value = $pair[1]
}

Write-Host "Columns added to sheets 4 and 5."
